$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = "002/DR002"
$ws.Range("I2").Value = 100000
$ws.Range("J2").Value = 7500
$ws.Range("K2").Value = 92500

# Row 3 updates
$ws.Range("A3").Value = "Karim benzima"
$ws.Range("B3").Value = "KS10293"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "119298479343984738347747"
$ws.Range("E3").Value = "BMCE"
$ws.Range("G3").Value = "794/DR KESH"
$ws.Range("H3").Value = "annuelle"
$ws.Range("I3").Value = 100000
$ws.Range("J3").Value = 10000
$ws.Range("K3").Value = 90000

# Row 4 (totals) updates
$ws.Range("I4").Value = 200000
$ws.Range("J4").Value = 17500
$ws.Range("K4").Value = 182500
